$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 427 (old rows 427-446
# shift down to 429-448), matching the new weekly price-report entries.
$ws.Rows.Item(427).Insert()
$ws.Rows.Item(427).Insert()

# New row 427
$ws.Cells.Item(427, 1).Value = 6
$ws.Cells.Item(427, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(427, 3).Value = "Metropolitana"
$ws.Cells.Item(427, 4).Value = 44714
$ws.Cells.Item(427, 5).Value = 13
$ws.Cells.Item(427, 6).Value = 100112043
$ws.Cells.Item(427, 7).Value = "Pepino ensalada"
$ws.Cells.Item(427, 8).Value = "Sin especificar"
$ws.Cells.Item(427, 9).Value = "Primera"
$ws.Cells.Item(427, 10).Value = 200
$ws.Cells.Item(427, 11).Value = 18000
$ws.Cells.Item(427, 12).Value = 20000
$ws.Cells.Item(427, 13).Value = 18800
$ws.Cells.Item(427, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(427, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(427, 16).Value = 313
$ws.Cells.Item(427, 17).Value = 60
$ws.Cells.Item(427, 18).Value = "Hortaliza"

# New row 428
$ws.Cells.Item(428, 1).Value = 6
$ws.Cells.Item(428, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(428, 3).Value = "Metropolitana"
$ws.Cells.Item(428, 4).Value = 44714
$ws.Cells.Item(428, 5).Value = 13
$ws.Cells.Item(428, 6).Value = 100112043
$ws.Cells.Item(428, 7).Value = "Pepino ensalada"
$ws.Cells.Item(428, 8).Value = "Sin especificar"
$ws.Cells.Item(428, 9).Value = "Segunda"
$ws.Cells.Item(428, 10).Value = 230
$ws.Cells.Item(428, 11).Value = 16000
$ws.Cells.Item(428, 12).Value = 18000
$ws.Cells.Item(428, 13).Value = 16870
$ws.Cells.Item(428, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(428, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(428, 16).Value = 211
$ws.Cells.Item(428, 17).Value = 80
$ws.Cells.Item(428, 18).Value = "Hortaliza"
